$d = $word.ActiveDocument

# 1) Objetivos body: replace with 'Programa resumido' text (para 6)
$p = $d.Paragraphs.Item(6)
$d.Range($p.Range.Start, $p.Range.End).Text = '1.Introduçâo 2. Conceitos básicos sobre materiais compósitos, suas matrizes e seus processo de fabricação 3. Tipos de reforços 4. Compósitos nanoestruturados, naturais e híbridos 5. Mecânica da estrutura reforçada 6. Atividade prática'

# 2) Docente(s) list paragraph: rebuild with Objetivos text + Programa text + Metodo + Criterio (para 8)
$p = $d.Paragraphs.Item(8)
$d.Range($p.Range.Start, $p.Range.End).Text = 'Fornecer aos estudantes uma visão abrangente e interdisciplinar sobre materiais compósitos, além de mostrar as especificidades de cada matriz, sendo ela metálica, cerâmica ou polimérica. Ademais, deseja-se apresentar os fundamentos teóricos da mecânica de estruturas reforçadas e a partir de atividades práticas demostrar métodos de caracterização de materiais compósitos e como prepara-los.' + [char]11 + '1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. ' + [char]11 + '2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. ' + [char]11 + '3. Conceitos de Interface' + [char]11 + '4. Compósitos de matriz metálica: características e processos de fabricação. ' + [char]11 + '5. Compósitos de matriz cerâmica: características e processos de fabricação. ' + [char]11 + '6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. ' + [char]11 + '7. Compósitos nanoestruturados. ' + [char]11 + '8. Compósitos Naturais. ' + [char]11 + '9. Compósitos Híbridos ' + [char]11 + '10. Mecânica de estruturas reforçadas. ' + [char]11 + 'Conteúdo prático: ' + [char]11 + '1. Caracterização e análise de compósitos de matriz metálica. ' + [char]11 + '2. Preparação e caracterização de compósitos de matriz polimérica.' + [char]11 + '(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) ' + [char]11 + '3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados' + [char]11 + 'De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.' + [char]11 + 'A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)'

# 3) 'Programa resumido' body: replace with Norma de recuperacao text (para 10)
$p = $d.Paragraphs.Item(10)
$d.Range($p.Range.Start, $p.Range.End).Text = 'Devido a cunho prático da disciplina não haverá recuperação.'

# 4) 'Programa' body: replace with Bibliografia text (para 12)
$p = $d.Paragraphs.Item(12)
$d.Range($p.Range.Start, $p.Range.End).Text = '1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Compósitos estruturais: tecnologia e prática. São Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.'

# 5) Avaliacao body: replace the three non-bold value runs with the docente names (para 14)
$p14 = $d.Paragraphs.Item(14)
$scope = $p14.Range
$scope.Find.Execute('De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.', $false, $false, $false, $false, $false, $true, 1, $false, '519033 - Carlos Yujiro Shigue', 2) | Out-Null
$scope = $p14.Range
$scope.Find.Execute('A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)', $false, $false, $false, $false, $false, $true, 1, $false, '3586455 - Cassius Olivio Figueiredo Terra Ruchert', 2) | Out-Null
$scope = $p14.Range
$scope.Find.Execute('Devido a cunho prático da disciplina não haverá recuperação.', $false, $false, $false, $false, $false, $true, 1, $false, '1033242 - Fábio Herbst Florenzano', 2) | Out-Null

# 6) Bibliografia body: replace with last docente name (para 16)
$p = $d.Paragraphs.Item(16)
$d.Range($p.Range.Start, $p.Range.End).Text = '1922320 - Sebastiao Ribeiro'

